# Backup QR Scanner data - append new scan log rows and rename the sheet
# to match the subject being logged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab from "Scanner" to "Pathology_Lab_Museum"
$ws.Name = "Pathology_Lab_Museum"

# Append the two new QR-scan log entries right after the current last row (127)
$newRows = @(
    @("244030", "Pathology Lab/Museum", "18/11/2025", "09:38:57", "Scan", "mona.I.hussein@med.asu.edu.eg"),
    @("244047", "Pathology Lab/Museum", "18/11/2025", "09:39:03", "Scan", "mona.I.hussein@med.asu.edu.eg")
)

$lastRow = 127
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $lastRow + 1 + $i
    $row = $newRows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($c -eq 1) {
            # Student ID column holds numeric-looking text, not a real
            # number, matching the rest of the column (e.g. A2:A127).
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$c - 1]
    }
}
